$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.998.55"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.467.14"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.52"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.37"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "2.905.03"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "57.924.74"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.31"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "2.468.34"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.77"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.73"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.08"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.411"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.33"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.27"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  -4.72%  "
$ws.Range("E31").Value = "  -6.42%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.00"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("E36").Value = "  -10.37%  "
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("E38").Value = "  -4.79%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "272.78"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.00"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.68"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.01"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").Value = "1.729.49"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E51").Value = "  -1.11%  "
